$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set value for A1
$ws.Range("A1").Value = "Band Name"

# Copy formatting from B1 to A1 (single-cell copy, safe)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match column A's width to column B (same custom width, same column-level style)
$ws.Range("A1").EntireColumn.ColumnWidth = $ws.Range("B1").EntireColumn.ColumnWidth
